# Auto-generated COM-interop script to apply the 2022-Q1 sheet addition
$wb = $excel.ActiveWorkbook

# --- Step 1: repurpose the existing last sheet ('总计') to become '2022-Q1' ---
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = '2022-Q1'

# --- Step 2: add a brand-new sheet after it, named '总计' (the new running total sheet) ---
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$totalSheet.Name = '总计'


function Set-HeaderCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-IndexCell($ws, $row, $col, $num) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $num
    $cell.Font.Bold = $true
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-NumberCell($ws, $row, $col, $num) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $num
}

# --- Step 3: populate the '2022-Q1' sheet (fund holdings detail) ---
$ws1 = $q1
Set-HeaderCell $ws1 1 2 '基金代码'
Set-HeaderCell $ws1 1 3 '基金名称'
Set-HeaderCell $ws1 1 4 '基金规模'
Set-HeaderCell $ws1 1 5 '股票总仓位'
Set-HeaderCell $ws1 1 6 '仓位占比'
Set-HeaderCell $ws1 1 7 '持有市值(亿元)'
Set-HeaderCell $ws1 1 8 '仓位排名'

Set-IndexCell $ws1 2 1 0
Set-TextCell $ws1 2 2 '011300'
Set-TextCell $ws1 2 3 '易方达智造优势混合A'
Set-TextCell $ws1 2 4 '63.40'
Set-TextCell $ws1 2 5 '89.88'
Set-TextCell $ws1 2 6 '5.85'
Set-TextCell $ws1 2 7 '3.7089'
Set-NumberCell $ws1 2 8 6

Set-IndexCell $ws1 3 1 1
Set-TextCell $ws1 3 2 '011822'
Set-TextCell $ws1 3 3 '易方达产业升级一年封闭运作混合型证券投资基金A'
Set-TextCell $ws1 3 4 '65.20'
Set-TextCell $ws1 3 5 '64.75'
Set-TextCell $ws1 3 6 '3.22'
Set-TextCell $ws1 3 7 '2.0994'
Set-NumberCell $ws1 3 8 9

Set-IndexCell $ws1 4 1 2
Set-TextCell $ws1 4 2 '001182'
Set-TextCell $ws1 4 3 '易方达安心回馈混合'
Set-TextCell $ws1 4 4 '90.36'
Set-TextCell $ws1 4 5 '34.92'
Set-TextCell $ws1 4 6 '1.86'
Set-TextCell $ws1 4 7 '1.6807'
Set-NumberCell $ws1 4 8 3

Set-IndexCell $ws1 5 1 3
Set-TextCell $ws1 5 2 '007130'
Set-TextCell $ws1 5 3 '中庚小盘价值股票'
Set-TextCell $ws1 5 4 '40.99'
Set-TextCell $ws1 5 5 '93.10'
Set-TextCell $ws1 5 6 '3.02'
Set-TextCell $ws1 5 7 '1.2379'
Set-NumberCell $ws1 5 8 4

Set-IndexCell $ws1 6 1 4
Set-TextCell $ws1 6 2 '001373'
Set-TextCell $ws1 6 3 '易方达新丝路灵活配置混合'
Set-TextCell $ws1 6 4 '38.95'
Set-TextCell $ws1 6 5 '90.79'
Set-TextCell $ws1 6 6 '2.94'
Set-TextCell $ws1 6 7 '1.1451'
Set-NumberCell $ws1 6 8 10

Set-IndexCell $ws1 7 1 5
Set-TextCell $ws1 7 2 '512400'
Set-TextCell $ws1 7 3 '南方中证申万有色金属ETF'
Set-TextCell $ws1 7 4 '36.45'
Set-TextCell $ws1 7 5 '99.71'
Set-TextCell $ws1 7 6 '2.23'
Set-TextCell $ws1 7 7 '0.8128'
Set-NumberCell $ws1 7 8 10

Set-IndexCell $ws1 8 1 6
Set-TextCell $ws1 8 2 '012719'
Set-TextCell $ws1 8 3 '华夏新兴经济一年持有期混合型证券投资基金A'
Set-TextCell $ws1 8 4 '21.62'
Set-TextCell $ws1 8 5 '91.19'
Set-TextCell $ws1 8 6 '3.72'
Set-TextCell $ws1 8 7 '0.8043'
Set-NumberCell $ws1 8 8 2

Set-IndexCell $ws1 9 1 7
Set-TextCell $ws1 9 2 '011301'
Set-TextCell $ws1 9 3 '易方达智造优势混合C'
Set-TextCell $ws1 9 4 '9.96'
Set-TextCell $ws1 9 5 '89.88'
Set-TextCell $ws1 9 6 '5.85'
Set-TextCell $ws1 9 7 '0.5827'
Set-NumberCell $ws1 9 8 6

Set-IndexCell $ws1 10 1 8
Set-TextCell $ws1 10 2 '160221'
Set-TextCell $ws1 10 3 '国泰国证有色金属行业指数（LOF）A'
Set-TextCell $ws1 10 4 '25.54'
Set-TextCell $ws1 10 5 '94.76'
Set-TextCell $ws1 10 6 '2.28'
Set-TextCell $ws1 10 7 '0.5823'
Set-NumberCell $ws1 10 8 9

Set-IndexCell $ws1 11 1 9
Set-TextCell $ws1 11 2 '003624'
Set-TextCell $ws1 11 3 '创金合信资源主题精选股票A'
Set-TextCell $ws1 11 4 '9.51'
Set-TextCell $ws1 11 5 '89.59'
Set-TextCell $ws1 11 6 '3.68'
Set-TextCell $ws1 11 7 '0.3500'
Set-NumberCell $ws1 11 8 10

Set-IndexCell $ws1 12 1 10
Set-TextCell $ws1 12 2 '011823'
Set-TextCell $ws1 12 3 '易方达产业升级一年封闭运作混合型证券投资基金C'
Set-TextCell $ws1 12 4 '10.22'
Set-TextCell $ws1 12 5 '64.75'
Set-TextCell $ws1 12 6 '3.22'
Set-TextCell $ws1 12 7 '0.3291'
Set-NumberCell $ws1 12 8 9

Set-IndexCell $ws1 13 1 11
Set-TextCell $ws1 13 2 '003625'
Set-TextCell $ws1 13 3 '创金合信资源主题精选股票C'
Set-TextCell $ws1 13 4 '4.04'
Set-TextCell $ws1 13 5 '89.59'
Set-TextCell $ws1 13 6 '3.68'
Set-TextCell $ws1 13 7 '0.1487'
Set-NumberCell $ws1 13 8 10

Set-IndexCell $ws1 14 1 12
Set-TextCell $ws1 14 2 '011630'
Set-TextCell $ws1 14 3 '西藏东财中证有色金属指数增强A'
Set-TextCell $ws1 14 4 '3.87'
Set-TextCell $ws1 14 5 '92.28'
Set-TextCell $ws1 14 6 '3.49'
Set-TextCell $ws1 14 7 '0.1351'
Set-NumberCell $ws1 14 8 8

Set-IndexCell $ws1 15 1 13
Set-TextCell $ws1 15 2 '001959'
Set-TextCell $ws1 15 3 '华商乐享互联灵活配置混合'
Set-TextCell $ws1 15 4 '3.74'
Set-TextCell $ws1 15 5 '88.84'
Set-TextCell $ws1 15 6 '3.46'
Set-TextCell $ws1 15 7 '0.1294'
Set-NumberCell $ws1 15 8 3

Set-IndexCell $ws1 16 1 14
Set-TextCell $ws1 16 2 '003839'
Set-TextCell $ws1 16 3 '易方达瑞通灵活配置混合A'
Set-TextCell $ws1 16 4 '9.33'
Set-TextCell $ws1 16 5 '32.04'
Set-TextCell $ws1 16 6 '1.27'
Set-TextCell $ws1 16 7 '0.1185'
Set-NumberCell $ws1 16 8 8

Set-IndexCell $ws1 17 1 15
Set-TextCell $ws1 17 2 '003882'
Set-TextCell $ws1 17 3 '易方达瑞弘灵活配置混合A'
Set-TextCell $ws1 17 4 '6.90'
Set-TextCell $ws1 17 5 '28.58'
Set-TextCell $ws1 17 6 '1.56'
Set-TextCell $ws1 17 7 '0.1076'
Set-NumberCell $ws1 17 8 3

Set-IndexCell $ws1 18 1 16
Set-TextCell $ws1 18 2 '001891'
Set-TextCell $ws1 18 3 '中欧成长优选回报灵活配置混合E'
Set-TextCell $ws1 18 4 '2.97'
Set-TextCell $ws1 18 5 '94.42'
Set-TextCell $ws1 18 6 '2.57'
Set-TextCell $ws1 18 7 '0.0763'
Set-NumberCell $ws1 18 8 9

Set-IndexCell $ws1 19 1 17
Set-TextCell $ws1 19 2 '166020'
Set-TextCell $ws1 19 3 '中欧成长优选回报灵活配置混合A'
Set-TextCell $ws1 19 4 '2.97'
Set-TextCell $ws1 19 5 '94.42'
Set-TextCell $ws1 19 6 '2.57'
Set-TextCell $ws1 19 7 '0.0763'
Set-NumberCell $ws1 19 8 9

Set-IndexCell $ws1 20 1 18
Set-TextCell $ws1 20 2 '011631'
Set-TextCell $ws1 20 3 '西藏东财中证有色金属指数增强C'
Set-TextCell $ws1 20 4 '2.18'
Set-TextCell $ws1 20 5 '92.28'
Set-TextCell $ws1 20 6 '3.49'
Set-TextCell $ws1 20 7 '0.0761'
Set-NumberCell $ws1 20 8 8

Set-IndexCell $ws1 21 1 19
Set-TextCell $ws1 21 2 '012720'
Set-TextCell $ws1 21 3 '华夏新兴经济一年持有期混合型证券投资基金C'
Set-TextCell $ws1 21 4 '1.04'
Set-TextCell $ws1 21 5 '91.19'
Set-TextCell $ws1 21 6 '3.72'
Set-TextCell $ws1 21 7 '0.0387'
Set-NumberCell $ws1 21 8 2

Set-IndexCell $ws1 22 1 20
Set-TextCell $ws1 22 2 '003883'
Set-TextCell $ws1 22 3 '易方达瑞弘灵活配置混合C'
Set-TextCell $ws1 22 4 '2.02'
Set-TextCell $ws1 22 5 '28.58'
Set-TextCell $ws1 22 6 '1.56'
Set-TextCell $ws1 22 7 '0.0315'
Set-NumberCell $ws1 22 8 3

Set-IndexCell $ws1 23 1 21
Set-TextCell $ws1 23 2 '159881'
Set-TextCell $ws1 23 3 '国泰中证有色金属交易型开放式指数证券投资基金'
Set-TextCell $ws1 23 4 '1.29'
Set-TextCell $ws1 23 5 '98.16'
Set-TextCell $ws1 23 6 '2.08'
Set-TextCell $ws1 23 7 '0.0268'
Set-NumberCell $ws1 23 8 10

Set-IndexCell $ws1 24 1 22
Set-TextCell $ws1 24 2 '001266'
Set-TextCell $ws1 24 3 '国投瑞银招财灵活配置混合'
Set-TextCell $ws1 24 4 '0.76'
Set-TextCell $ws1 24 5 '67.47'
Set-TextCell $ws1 24 6 '2.34'
Set-TextCell $ws1 24 7 '0.0178'
Set-NumberCell $ws1 24 8 8

Set-IndexCell $ws1 25 1 23
Set-TextCell $ws1 25 2 '003840'
Set-TextCell $ws1 25 3 '易方达瑞通灵活配置混合C'
Set-TextCell $ws1 25 4 '1.38'
Set-TextCell $ws1 25 5 '32.04'
Set-TextCell $ws1 25 6 '1.27'
Set-TextCell $ws1 25 7 '0.0175'
Set-NumberCell $ws1 25 8 8

Set-IndexCell $ws1 26 1 24
Set-TextCell $ws1 26 2 '159876'
Set-TextCell $ws1 26 3 '华宝中证有色金属ETF'
Set-TextCell $ws1 26 4 '0.45'
Set-TextCell $ws1 26 5 '98.90'
Set-TextCell $ws1 26 6 '2.08'
Set-TextCell $ws1 26 7 '0.0094'
Set-NumberCell $ws1 26 8 10

Set-IndexCell $ws1 27 1 25
Set-TextCell $ws1 27 2 '159871'
Set-TextCell $ws1 27 3 '银华中证有色金属ETF'
Set-TextCell $ws1 27 4 '0.43'
Set-TextCell $ws1 27 5 '97.56'
Set-TextCell $ws1 27 6 '2.07'
Set-TextCell $ws1 27 7 '0.0089'
Set-NumberCell $ws1 27 8 10

Set-IndexCell $ws1 28 1 26
Set-TextCell $ws1 28 2 '159880'
Set-TextCell $ws1 28 3 '鹏华国证有色金属行业ETF'
Set-TextCell $ws1 28 4 '0.33'
Set-TextCell $ws1 28 5 '96.25'
Set-TextCell $ws1 28 6 '2.34'
Set-TextCell $ws1 28 7 '0.0077'
Set-NumberCell $ws1 28 8 9

Set-IndexCell $ws1 29 1 27
Set-TextCell $ws1 29 2 '012297'
Set-TextCell $ws1 29 3 '东兴宸瑞量化混合型证券投资基金A'
Set-TextCell $ws1 29 4 '0.48'
Set-TextCell $ws1 29 5 '89.42'
Set-TextCell $ws1 29 6 '1.57'
Set-TextCell $ws1 29 7 '0.0075'
Set-NumberCell $ws1 29 8 4

Set-IndexCell $ws1 30 1 28
Set-TextCell $ws1 30 2 '516650'
Set-TextCell $ws1 30 3 '华夏中证细分有色金属产业主题交易型开放式指数证券投资基金'
Set-TextCell $ws1 30 4 '0.29'
Set-TextCell $ws1 30 5 '98.99'
Set-TextCell $ws1 30 6 '2.20'
Set-TextCell $ws1 30 7 '0.0064'
Set-NumberCell $ws1 30 8 10

Set-IndexCell $ws1 31 1 29
Set-TextCell $ws1 31 2 '012298'
Set-TextCell $ws1 31 3 '东兴宸瑞量化混合型证券投资基金C'
Set-TextCell $ws1 31 4 '0.21'
Set-TextCell $ws1 31 5 '89.42'
Set-TextCell $ws1 31 6 '1.57'
Set-TextCell $ws1 31 7 '0.0033'
Set-NumberCell $ws1 31 8 4

# --- Step 4: populate the '总计' (Total) sheet ---
$ws2 = $totalSheet
Set-HeaderCell $ws2 1 2 '日期'
Set-HeaderCell $ws2 1 3 '持有数量(只)'
Set-HeaderCell $ws2 1 4 '持有市值(亿元)'

Set-IndexCell $ws2 2 1 0
Set-TextCell $ws2 2 2 '2022-Q1'
Set-NumberCell $ws2 2 3 30
Set-NumberCell $ws2 2 4 14.38

Set-IndexCell $ws2 3 1 1
Set-TextCell $ws2 3 2 '2021-Q4'
Set-NumberCell $ws2 3 3 52
Set-NumberCell $ws2 3 4 27.8

Set-IndexCell $ws2 4 1 2
Set-TextCell $ws2 4 2 '2021-Q3'
Set-NumberCell $ws2 4 3 37
Set-NumberCell $ws2 4 4 27.17

Set-IndexCell $ws2 5 1 3
Set-TextCell $ws2 5 2 '2021-Q2'
Set-NumberCell $ws2 5 3 28
Set-NumberCell $ws2 5 4 11.32

Set-IndexCell $ws2 6 1 4
Set-TextCell $ws2 6 2 '2021-Q1'
Set-NumberCell $ws2 6 3 20
Set-NumberCell $ws2 6 4 5.94

Set-IndexCell $ws2 7 1 5
Set-TextCell $ws2 7 2 '2020-Q4'
Set-NumberCell $ws2 7 3 21
Set-NumberCell $ws2 7 4 5.01

Write-Output "done"